$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New archival record (MCH238) being appended as row 2, right under the
# header row. Build the row's font (Calibri 10pt, theme text color) once on
# A2, then fan it out to the rest of the row via format paint so every cell
# ends up sharing the same style.
$a2 = $ws.Range("A2")
$a2.Font.ThemeColor = 1
$a2.Font.Name = "Calibri"
$a2.Font.Size = 10

$a2.Copy() | Out-Null
$ws.Range("C2:H2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Fill in the cell values. D2 and H2 are intentionally left blank (they are
# spacer columns) but keep the row's new style from the paste above.
$ws.Range("A2").Value = "MCH238"
$ws.Range("C2").Value = "CERAMICS  AS HISTORY OR CERAMICS AS OBJECTS"
$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: CABINET 1C | GRAP COUNT NUMER: NONE"

# Column B has no data in this row -- drop it entirely so it doesn't linger
# as an empty, styled cell.
$ws.Range("B2").Clear() | Out-Null

# Mirror the selection left behind after entering the new row of data.
$ws.Range("A2:K2").Select() | Out-Null
